$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '69.279.86'
$cell.Style = $origStyle

$ws.Range("E2").Value = '  +1.50%  '
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.388.24'
$cell.Style = $origStyle

$ws.Range("E3").Value = '  +1.23%  '
$ws.Range("E4").Value = '  +0.07%  '
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '581.09'
$cell.Style = $origStyle

$ws.Range("E5").Value = '  -0.33%  '
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '179.17'
$cell.Style = $origStyle

$ws.Range("E6").Value = '  +1.18%  '
$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = $origStyle

$ws.Range("E7").Value = '  +0.00%  '
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.593'
$cell.Style = $origStyle

$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("E9").Value = '  +8.32%  '
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.586'
$cell.Style = $origStyle

$ws.Range("E10").Value = '  +0.89%  '
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '48.35'
$cell.Style = $origStyle

$ws.Range("E11").Value = '  +0.90%  '
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0000282'
$cell.Style = $origStyle

$ws.Range("E12").Value = '  +3.34%  '
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '683.38'
$cell.Style = $origStyle

$ws.Range("E13").Value = '  -0.07%  '
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.58'
$cell.Style = $origStyle

$ws.Range("E14").Value = '  +2.05%  '
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.932.20'
$cell.Style = $origStyle

$ws.Range("E15").Value = '  +0.97%  '
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '69.409.97'
$cell.Style = $origStyle

$ws.Range("E16").Value = '  +1.59%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.389.56'
$cell.Style = $origStyle

$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.120'
$cell.Style = $origStyle

$ws.Range("E18").Value = '  +0.85%  '
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.66'
$cell.Style = $origStyle

$ws.Range("E19").Value = '  +1.40%  '
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.28'
$cell.Style = $origStyle

$ws.Range("E20").Value = '  +0.96%  '
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.909'
$cell.Style = $origStyle

$ws.Range("E21").Value = '  +1.77%  '
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.13'
$cell.Style = $origStyle

$ws.Range("E22").Value = '  +1.38%  '
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.37'
$cell.Style = $origStyle

$ws.Range("E23").Value = '  -1.22%  '
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '101.09'
$cell.Style = $origStyle

$ws.Range("E24").Value = '  +1.27%  '
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("E26").Value = '  +0.38%  '
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.72'
$cell.Style = $origStyle

$ws.Range("E27").Value = '  +2.22%  '
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '33.41'
$cell.Style = $origStyle

$ws.Range("E28").Value = '  +1.44%  '
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.72'
$cell.Style = $origStyle

$ws.Range("E29").Value = '  +2.78%  '
$ws.Range("E30").Value = '  +0.42%  '
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.84'
$cell.Style = $origStyle

$ws.Range("E31").Value = '  +16.99%  '
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.03'
$cell.Style = $origStyle

$ws.Range("E32").Value = '  -0.19%  '
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '549.58'
$cell.Style = $origStyle

$ws.Range("E33").Value = '  -2.24%  '
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.105'
$cell.Style = $origStyle

$ws.Range("E34").Value = '  +0.22%  '
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '57.84'
$cell.Style = $origStyle

$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  +0.11%  '
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.602.65'
$cell.Style = $origStyle

$ws.Range("E37").Value = '  -2.68%  '
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.140'
$cell.Style = $origStyle

$ws.Range("E38").Value = '  +3.28%  '
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '35.42'
$cell.Style = $origStyle

$ws.Range("E39").Value = '  +2.26%  '
$ws.Range("E40").Value = '  +10.76%  '
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.31'
$cell.Style = $origStyle

$ws.Range("E41").Value = '  +5.03%  '
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.71'
$cell.Style = $origStyle

$ws.Range("E42").Value = '  +4.14%  '
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.37'
$cell.Style = $origStyle

$ws.Range("E43").Value = '  +3.62%  '
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0424'
$cell.Style = $origStyle

$ws.Range("E44").Value = '  +3.60%  '
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.335'
$cell.Style = $origStyle

$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("E46").Value = '  +1.11%  '
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("E48").Value = '  +3.68%  '
$ws.Range("E49").Value = '  -0.10%  '
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '129.83'
$cell.Style = $origStyle

$ws.Range("E50").Value = '  -1.08%  '
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.59'
$cell.Style = $origStyle

$ws.Range("E51").Value = '  +1.43%  '
